$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 353.8846
$ws.Range("I33").Value = 338.42856
$ws.Range("J33").Value = 371.91666
$ws.Range("K33").Value = 338.42856
$ws.Range("L33").Value = 371.91666
$ws.Range("M33").Value = -109.42856
$ws.Range("N33").Value = -829.91666
$ws.Range("H64").Value = 3740
$ws.Range("I64").Value = 4166.6665
$ws.Range("J64").Value = 3228
$ws.Range("K64").Value = 4166.6665
$ws.Range("L64").Value = 3228
$ws.Range("M64").Value = -3918.6665
$ws.Range("N64").Value = -3724
$ws.Range("H67").Value = 3740
$ws.Range("I67").Value = 4166.6665
$ws.Range("J67").Value = 3228
$ws.Range("K67").Value = 4166.6665
$ws.Range("L67").Value = 3228
$ws.Range("M67").Value = -3308.6665
$ws.Range("N67").Value = -4944
$ws.Range("H111").Value = 661.73334
$ws.Range("I111").Value = 562.7778
$ws.Range("J111").Value = 810.1667
$ws.Range("K111").Value = 1688.3334
$ws.Range("L111").Value = 2430.5001
$ws.Range("M111").Value = 1378.6666
$ws.Range("N111").Value = -8564.500100000001
$ws.Range("H135").Value = 854.04346
$ws.Range("I135").Value = 863.9524
$ws.Range("K135").Value = 7775.5716
$ws.Range("M135").Value = -5240.5716
$ws.Range("H138").Value = 1715.2858
$ws.Range("I138").Value = 999.7451
$ws.Range("J138").Value = 3635.9473
$ws.Range("K138").Value = 2999.2353
$ws.Range("L138").Value = 10907.8419
$ws.Range("M138").Value = 2140.7647
$ws.Range("N138").Value = -21187.8419
$ws.Range("H141").Value = 747488.25
$ws.Range("J141").Value = 1193562.8
$ws.Range("L141").Value = 3580688.4
$ws.Range("N141").Value = -3591048.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 30900
$ws.Range("J7").Value = 30900
$ws.Range("L7").Value = 30900
$ws.Range("N7").Value = -31128
$ws.Range("H32").Value = 3821.24
$ws.Range("I32").Value = 3411.5056
$ws.Range("J32").Value = 7136.364
$ws.Range("K32").Value = 3411.5056
$ws.Range("L32").Value = 7136.364
$ws.Range("M32").Value = -3124.5056
$ws.Range("N32").Value = -7710.364
$ws.Range("H64").Value = 29181.818
$ws.Range("J64").Value = 29181.818
$ws.Range("L64").Value = 29181.818
$ws.Range("N64").Value = -29677.818
$ws.Range("H67").Value = 29181.818
$ws.Range("J67").Value = 29181.818
$ws.Range("L67").Value = 29181.818
$ws.Range("N67").Value = -30897.818
$ws.Range("H122").Value = 1491.081
$ws.Range("I122").Value = 1110
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 3330
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = -880
$ws.Range("N122").Value = -12460
$ws.Range("H132").Value = 1623.9036
$ws.Range("I132").Value = 1438.8108
$ws.Range("J132").Value = 3145.7778
$ws.Range("K132").Value = 4316.4324
$ws.Range("L132").Value = 9437.3334
$ws.Range("M132").Value = -1786.4324
$ws.Range("N132").Value = -14497.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1311.3235
$ws.Range("I105").Value = 1273.76
$ws.Range("J105").Value = 1415.6666
$ws.Range("K105").Value = 1273.76
$ws.Range("L105").Value = 1415.6666
$ws.Range("M105").Value = 473.24
$ws.Range("N105").Value = -4909.6666
$ws.Range("H107").Value = 676
$ws.Range("I107").Value = 584.94116
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 584.94116
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = 1335.05884
$ws.Range("N107").Value = -5290
$ws.Range("H134").Value = 1968.4193
$ws.Range("I134").Value = 1724.862
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 5174.586
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -2639.586
$ws.Range("N134").Value = -21570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1734.3529
$ws.Range("I16").Value = 1542.7778
$ws.Range("J16").Value = 1949.875
$ws.Range("K16").Value = 1542.7778
$ws.Range("L16").Value = 1949.875
$ws.Range("M16").Value = -1255.7778
$ws.Range("N16").Value = -2523.875
$ws.Range("H113").Value = 1734.3529
$ws.Range("I113").Value = 1542.7778
$ws.Range("J113").Value = 1949.875
$ws.Range("K113").Value = 1542.7778
$ws.Range("L113").Value = 1949.875
$ws.Range("M113").Value = 627.2221999999999
$ws.Range("N113").Value = -6289.875
$ws.Range("H132").Value = 3041.6
$ws.Range("I132").Value = 3136
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 9408
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -6878
$ws.Range("N132").Value = -13760
$ws.Range("H133").Value = 17136.166
$ws.Range("J133").Value = 17136.166
$ws.Range("L133").Value = 17136.166
$ws.Range("N133").Value = -22196.166
$ws.Range("H134").Value = 6348.4116
$ws.Range("I134").Value = 7363.3076
$ws.Range("K134").Value = 22089.9228
$ws.Range("M134").Value = -19554.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1096.6666
$ws.Range("I15").Value = 145
$ws.Range("K15").Value = 435
$ws.Range("M15").Value = -295
$ws.Range("H47").Value = 1998.1428
$ws.Range("J47").Value = 3427
$ws.Range("L47").Value = 10281
$ws.Range("N47").Value = -11143
$ws.Range("H74").Value = 13590
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13590
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 40770
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -42892
$ws.Range("H77").Value = 13590
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13590
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 122310
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -132918

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 27781.7
$ws.Range("I102").Value = 1898.9474
$ws.Range("J102").Value = 51199.43
$ws.Range("K102").Value = 1898.9474
$ws.Range("L102").Value = 51199.43
$ws.Range("M102").Value = -276.9474
$ws.Range("N102").Value = -54443.43
$ws.Range("H104").Value = 33300
$ws.Range("J104").Value = 33300
$ws.Range("L104").Value = 33300
$ws.Range("H122").Value = 3734.1155
$ws.Range("I122").Value = 4059.1428
$ws.Range("J122").Value = 3614.3684
$ws.Range("K122").Value = 12177.4284
$ws.Range("L122").Value = 10843.1052
$ws.Range("M122").Value = -9727.428400000001
$ws.Range("N122").Value = -15743.1052
$ws.Range("H126").Value = 2550.1892
$ws.Range("I126").Value = 1390.8462
$ws.Range("J126").Value = 3178.1667
$ws.Range("K126").Value = 4172.5386
$ws.Range("L126").Value = 9534.500100000001
$ws.Range("M126").Value = -1702.5386
$ws.Range("N126").Value = -14474.5001
$ws.Range("H132").Value = 3806
$ws.Range("I132").Value = 4867.933
$ws.Range("K132").Value = 14603.799
$ws.Range("M132").Value = -12073.799

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1661.6957
$ws.Range("I7").Value = 1355.5454
$ws.Range("K7").Value = 1355.5454
$ws.Range("M7").Value = -1243.5454
$ws.Range("H124").Value = 27107.25
$ws.Range("J124").Value = 27107.25
$ws.Range("L124").Value = 27107.25
$ws.Range("N124").Value = -36927.25
$ws.Range("H126").Value = 1661.6957
$ws.Range("I126").Value = 1355.5454
$ws.Range("K126").Value = 4066.6362
$ws.Range("M126").Value = -1596.6362
$ws.Range("H132").Value = 5807.3237
$ws.Range("I132").Value = 2268
$ws.Range("J132").Value = 7998.3335
$ws.Range("K132").Value = 6804
$ws.Range("L132").Value = 23995.0005
$ws.Range("M132").Value = -4274
$ws.Range("N132").Value = -29055.0005
$ws.Range("H136").Value = 3038.56
$ws.Range("I136").Value = 2759.111
$ws.Range("J136").Value = 3757.1428
$ws.Range("K136").Value = 8277.332999999999
$ws.Range("L136").Value = 11271.4284
$ws.Range("M136").Value = -5727.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1736.25
$ws.Range("I81").Value = 1815
$ws.Range("K81").Value = 3630
$ws.Range("M81").Value = -2569
$ws.Range("H84").Value = 1736.25
$ws.Range("I84").Value = 1815
$ws.Range("K84").Value = 18150
$ws.Range("M84").Value = -12846
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 10000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 10000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -10812
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 10000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 10000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -12808
$ws.Range("H100").Value = 618
$ws.Range("I100").Value = 597.5
$ws.Range("K100").Value = 1195
$ws.Range("M100").Value = -654
$ws.Range("H132").Value = 30729.895
$ws.Range("I132").Value = 5305.4546
$ws.Range("J132").Value = 65688.5
$ws.Range("K132").Value = 15916.3638
$ws.Range("L132").Value = 197065.5
$ws.Range("M132").Value = -13386.3638
$ws.Range("N132").Value = -202125.5
